$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.423.47'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").Value = '2.426.70'
$ws.Range("E3").Value = '  -1.98%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '557.74'
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").Value = '160.07'
$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +6.36%  '

$ws.Range("E11").Value = '  -0.76%  '

$ws.Range("D12").Value = '4.63'
$ws.Range("E12").Value = '  -5.53%  '

$ws.Range("D13").Value = '68.288.88'
$ws.Range("E13").Value = '  -1.18%  '

$ws.Range("D14").Value = '2.871.11'
$ws.Range("E14").Value = '  -1.19%  '

$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.00'
$ws.Range("E16").Value = '  -3.22%  '

$ws.Range("D17").Value = '2.425.66'
$ws.Range("E17").Value = '  +0.28%  '

$ws.Range("E18").Value = '  -3.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '334.70'
$ws.Range("E19").Value = '  -1.47%  '

$ws.Range("E20").Value = '  -2.06%  '

$ws.Range("D21").Value = '3.82'
$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.90'
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").Value = '66.61'
$ws.Range("E24").Value = '  -1.09%  '

$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").Value = '2.549.07'
$ws.Range("E26").Value = '  -2.11%  '

$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").Value = '0.0₃0814'
$ws.Range("E28").Value = '  -1.47%  '

$ws.Range("E29").Value = '  -1.24%  '

$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").Value = '425.51'
$ws.Range("E31").Value = '  -1.90%  '

$ws.Range("E32").Value = '  -0.55%  '

$ws.Range("E33").Value = '  -1.65%  '

$ws.Range("D34").Value = '158.73'
$ws.Range("E34").Value = '  +0.66%  '

$ws.Range("D35").Value = '19.04'
$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("E37").Value = '  +0.20%  '

$ws.Range("D38").Value = '0.105'
$ws.Range("E38").Value = '  -4.76%  '

$ws.Range("E39").Value = '  -1.89%  '

$ws.Range("D40").Value = '4.32'
$ws.Range("E40").Value = '  -3.13%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("E42").Value = '  -1.67%  '

$ws.Range("D43").Value = '132.78'
$ws.Range("E43").Value = '  -0.32%  '

$ws.Range("D44").Value = '2.01'
$ws.Range("E44").Value = '  -3.65%  '

$ws.Range("E45").Value = '  -0.92%  '

$ws.Range("D46").Value = '0.0713'
$ws.Range("E46").Value = '  -0.58%  '

$ws.Range("E47").Value = '  -1.50%  '

$ws.Range("E48").Value = '  -1.65%  '

$ws.Range("D49").Value = '0.0913'
$ws.Range("E49").Value = '  -0.52%  '

$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("E51").Value = '  -2.27%  '

